$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new TC3 row (row 4) mirroring the structure of the existing rows
$ws.Range("A4").Value = "forgotPassword"
$ws.Range("B4").Value = "Admin"
$ws.Range("C4").Value = "admin123"
$ws.Range("D4").Value = "john"

# Update the active selection to D7 as a single cell
$ws.Range("D7").Select()
